# Add "test_createDMN_FuncKey" test case data as a new row (row 24) on the
# "PMTestData" sheet, including a new column D used for the DMN
# number_initiate command string.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PMTestData")

# New column D needs to be wide enough to show the DMN command text
# (matches the ~108-char width used for similar wrapped text columns).
$ws.Columns.Item(4).ColumnWidth = 107.66

$ws.Cells.Item(24, 1).Value = "test_createDMN_FuncKey"

$ws.Cells.Item(24, 4).Value = "number_initiate -number 60000 -numbertype ex,`nextension -i -d 60000 -l 1 --csp 0,`nip_extension -i -d 60000"
$ws.Cells.Item(24, 4).WrapText = $true

$ws.Cells.Item(24, 2).Value = "number_initiate -number 90000 -numbertype ex,90000,90000,60000,1,FirstName,LastName,Mitel 6869i,DMN,ip_extension -e -d 60000..90000,extension -e -d 60000..90000,number_end -number 60000 -numbertype ex,number_end -number 90000 -numbertype ex"
$ws.Cells.Item(24, 2).WrapText = $true

$ws.Cells.Item(24, 3).Value = "Y"

$ws.Rows.Item(24).RowHeight = 72.5

$ws.Range("B24").Select()

$wb.Save()
